# Generate Report for Handback
#
# The handback for both locales (zh-cn, de-de) is now in sync with en-US:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Handback DateTime" is refreshed to the handback-generation time
#   - The stale "handback file is not the latest" warning in Error Detail is cleared
#   - Column widths are refreshed (auto-fit) to comfortably fit the new status text
#     and the now-empty error column.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de summary columns (E2, F2) mirror each locale's
# Status column.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.084
$overview.Columns.Item(6).ColumnWidth = 29.084

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-29 00:48:13"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.084
$zhcn.Columns.Item(16).ColumnWidth = 12.7501

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-29 00:48:20"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.084
$dede.Columns.Item(16).ColumnWidth = 12.7501
